# "delirio místico de signos" — flip the sign on a handful of derived
# measurement formulas on Hoja1 (S11, U11, I28, K28, I33, K33). Their
# dependents (W11, M28, M33, and the summary K43/K44/K46/K47/K49/K50)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("S11").Formula = "=-B11/E11*1000"
$ws.Range("U11").Formula = "=-D11/E11"

$ws.Range("I28").Formula = "=-B27/E27*1000"
$ws.Range("K28").Formula = "=-D27/E27"

$ws.Range("I33").Formula = "=-B32/E32*1000"
$ws.Range("K33").Formula = "=-D32/E32"

# Restore the view: scroll position and active selection.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("K34").Select()
